{"js": "// Replace the date line and each \"AxB=\" multiplication-problem cell\n// text with its updated value. Every <w:t> run in the document changes\n// exactly once, and all old/new values are unique, so a direct\n// search-and-replace per pair is unambiguous.\nconst replacements = [\n  [\"2024-05-08 Wednesday\", \"2024-05-09 Thursday\"],\n  [\"922\u00d77=\", \"962\u00d79=\"],\n  [\"664\u00d74=\", \"609\u00d73=\"],\n  [\"727\u00d78=\", \"544\u00d78=\"],\n  [\"665\u00d76=\", \"605\u00d79=\"],\n  [\"405\u00d79=\", \"994\u00d73=\"],\n  [\"323\u00d76=\", \"983\u00d76=\"],\n  [\"132\u00d75=\", \"400\u00d76=\"],\n  [\"498\u00d74=\", \"859\u00d75=\"],\n  [\"983\u00d72=\", \"275\u00d79=\"],\n  [\"318\u00d73=\", \"929\u00d79=\"],\n  [\"999\u00d74=\", \"594\u00d73=\"],\n  [\"231\u00d74=\", \"217\u00d74=\"],\n  [\"663\u00d73=\", \"863\u00d76=\"],\n  [\"357\u00d72=\", \"340\u00d72=\"],\n  [\"638\u00d74=\", \"145\u00d78=\"],\n  [\"135\u00d77=\", \"453\u00d77=\"],\n  [\"319\u00d79=\", \"919\u00d75=\"],\n  [\"326\u00d78=\", \"132\u00d73=\"],\n  [\"254\u00d77=\", \"215\u00d79=\"],\n  [\"647\u00d76=\", \"294\u00d73=\"],\n  [\"572\u00d76=\", \"866\u00d73=\"],\n  [\"359\u00d79=\", \"542\u00d73=\"],\n  [\"556\u00d77=\", \"658\u00d78=\"],\n  [\"562\u00d72=\", \"321\u00d76=\"],\n  [\"333\u00d79=\", \"667\u00d79=\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n", "ps1": "# Replace the date line and each \"AxB=\" multiplication-problem cell\n# text with its updated value. Every run in the document changes exactly\n# once, and all old/new values are unique, so Find/Replace per pair is\n# unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-05-08 Wednesday\", \"2024-05-09 Thursday\"),\n    @(\"922\u00d77=\", \"962\u00d79=\"),\n    @(\"664\u00d74=\", \"609\u00d73=\"),\n    @(\"727\u00d78=\", \"544\u00d78=\"),\n    @(\"665\u00d76=\", \"605\u00d79=\"),\n    @(\"405\u00d79=\", \"994\u00d73=\"),\n    @(\"323\u00d76=\", \"983\u00d76=\"),\n    @(\"132\u00d75=\", \"400\u00d76=\"),\n    @(\"498\u00d74=\", \"859\u00d75=\"),\n    @(\"983\u00d72=\", \"275\u00d79=\"),\n    @(\"318\u00d73=\", \"929\u00d79=\"),\n    @(\"999\u00d74=\", \"594\u00d73=\"),\n    @(\"231\u00d74=\", \"217\u00d74=\"),\n    @(\"663\u00d73=\", \"863\u00d76=\"),\n    @(\"357\u00d72=\", \"340\u00d72=\"),\n    @(\"638\u00d74=\", \"145\u00d78=\"),\n    @(\"135\u00d77=\", \"453\u00d77=\"),\n    @(\"319\u00d79=\", \"919\u00d75=\"),\n    @(\"326\u00d78=\", \"132\u00d73=\"),\n    @(\"254\u00d77=\", \"215\u00d79=\"),\n    @(\"647\u00d76=\", \"294\u00d73=\"),\n    @(\"572\u00d76=\", \"866\u00d73=\"),\n    @(\"359\u00d79=\", \"542\u00d73=\"),\n    @(\"556\u00d77=\", \"658\u00d78=\"),\n    @(\"562\u00d72=\", \"321\u00d76=\"),\n    @(\"333\u00d79=\", \"667\u00d79=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Could not find text to replace: $oldText\"\n    }\n}\n"}
